$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.208.00"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "'2.447.51"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'580.63"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("D6").Value = "'143.45"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.531"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "'2.444.19"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("D12").Value = "'5.20"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "'0.345"
$ws.Range("E13").Value = "  -2.28%  "
$ws.Range("D14").Value = "'26.35"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").Value = "'0.0000173"
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("D16").Value = "'2.800.58"
$ws.Range("E16").Value = "  -2.51%  "
$ws.Range("D17").Value = "'62.038.20"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "'2.421.49"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("D21").Value = "'328.67"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("D22").Value = "'4.10"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "'65.68"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("E26").Value = "  +7.41%  "
$ws.Range("D27").Value = "'607.98"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D29").Value = "'0.0₃0950"
$ws.Range("E29").Value = "  -5.13%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("E31").Value = "  -4.16%  "
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").Value = "'0.377"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("E38").Value = "  -4.80%  "
$ws.Range("D39").Value = "'149.42"
$ws.Range("E39").Value = "  +2.98%  "
$ws.Range("D40").Value = "'5.32"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("D43").Value = "'42.67"
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -4.21%  "
$ws.Range("D46").Value = "'143.17"
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("D47").Value = "'3.63"
$ws.Range("E47").Value = "  -2.72%  "
$ws.Range("D48").Value = "'0.605"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").Value = "'0.0523"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").Value = "'19.45"
$ws.Range("E50").Value = "  -6.16%  "
$ws.Range("E51").Value = "  +9.70%  "
